$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to B,C,D,E,F,I,J,K,L,M,N (skipping G=const 1 and H=empty)
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @(
    @(1.02,1.03056934243225,1.039711758461382,0.992614727750844,1.046522837754324,1.036417295402131,1.035709862428591,1.042496168756128,0.9955398523335997,1.049288036417654,1.015910621668556),
    @(1.02,1.031314871927964,1.040307387316647,0.9936372048519299,1.04729102970642,1.036562709037499,1.036097941903592,1.042902638045458,0.9963617723202687,1.049868000960313,1.016039493011981),
    @(1.02,1.031797890861302,1.040693350401313,0.9942998659930998,1.047789133064154,1.036655910832322,1.036348970242741,1.043165503753027,0.9968940712668347,1.050243645339527,1.016122837147318),
    @(1.02,1.032001096481772,1.040855739320687,0.994578699834602,1.04799878053123,1.036694878916183,1.036454480898137,1.043275975848731,0.9971179600053012,1.050401652656046,1.016157863973619),
    @(1.02,1.032035224014672,1.040883012703132,0.994625531979634,1.04803399555557,1.036701409260948,1.036472195293524,1.043294522409025,0.9971555583673455,1.050428187794795,1.016163744470717),
    @(1.02,1.031800605539538,1.040695519741513,0.994303590798249,1.047791933424677,1.036656432367382,1.036350380168019,1.043166980032375,0.9968970624462089,1.050245756303728,1.016123305221472),
    @(1.02,1.030821169838965,1.039912938673287,0.9929600610674297,1.046782236859284,1.036466622563597,1.035841032528063,1.042633566553951,0.9958175282591056,1.049483960754287,1.015954183228435),
    @(1.02,1.029100047882665,1.038538240625641,0.9906006454969559,1.045011022084601,1.036125366151758,1.034942903233105,1.041692560507693,0.9939188001724441,1.048144482254137,1.01565585002003),
    @(1.02,1.02795595516783,1.037624790379215,0.989033133672735,1.04383572226333,1.035893343492465,1.03434382903707,1.04106459379581,0.9926553831429383,1.047253556363683,1.015456774259791),
    @(1.02,1.027461362771713,1.037229996403466,0.988355674866747,1.043328139279059,1.035791812943961,1.034084363562125,1.040792546138683,0.9921088820399291,1.046868287116006,1.015370533692701),
    @(1.02,1.027277772404786,1.037083464950837,0.9881042295826724,1.043139802391886,1.03575394100228,1.033987978491677,1.040691476873128,0.9919059725120875,1.046725259081031,1.015338494647061),
    @(1.02,1.027317147551598,1.037114891316513,0.9881581567098651,1.043180192114548,1.035762071844096,1.034008653758916,1.040713157394197,0.9919494934313052,1.046755935515134,1.015345367370648),
    @(1.02,1.027446184608554,1.037217881758977,0.9883348863814464,1.043312567158637,1.0357886856761,1.034076396498445,1.040784192093688,0.9920921077337197,1.046856462769658,1.015367885443857),
    @(1.02,1.027525704992002,1.037281352552051,0.9884438009545853,1.04339415463166,1.035805062293326,1.034118134010512,1.040827956488704,0.9921799884222134,1.046918411336105,1.015381758862536),
    @(1.02,1.027988796814391,1.037651007235495,0.9890781214508737,1.043869437081742,1.035900059404851,1.034361047705402,1.041082646006521,0.9926916645766087,1.047279136267389,1.015462496970985),
    @(1.02,1.028279499616524,1.037883080337851,0.989476357848556,1.044167926784074,1.035959364468392,1.034513405285543,1.04124237119311,0.9930127773699352,1.047505546687327,1.015513131622003),
    @(1.02,1.028449139402355,1.038018515530271,0.9897087662937556,1.044342158911403,1.035993853473656,1.03460226677478,1.04133552317172,0.9932001317071769,1.047637656900184,1.015542662119673),
    @(1.02,1.028506995314499,1.038064707391242,0.9897880325774034,1.044401589228704,1.036005595901423,1.034632565169339,1.041367283335988,0.9932640239640975,1.047682711296998,1.015552730591081),
    @(1.02,1.028248301899698,1.037858173744377,0.9894336180360679,1.044135888393976,1.035953012206875,1.034497059381418,1.041225235528388,0.9929783193494215,1.047481249926231,1.015507699396719),
    @(1.02,1.02740818299264,1.037187550515812,0.9882828385668249,1.043273580412133,1.0357808529483,1.034056448161345,1.040763274648305,0.9920501090198102,1.046826857827375,1.015361254579493),
    @(1.02,1.026880681080028,1.036766555593608,0.9875604150241495,1.042732583186613,1.035671690183964,1.033779372539856,1.040472714314323,0.9914670000341481,1.046415868429027,1.015269147556466),
    @(1.02,1.027160251465381,1.036989670438976,0.9879432794643023,1.043019264387198,1.035729646325931,1.033926259471149,1.040626755515774,0.991776070289318,1.04663369813449,1.01531797804813),
    @(1.02,1.028262398565041,1.037869427743152,0.9894529299347244,1.044150364767782,1.03595588283765,1.034504445410947,1.041232978438108,0.9929938892766442,1.047492228436056,1.015510153997418),
    @(1.02,1.029544422117994,1.038893110173925,0.9912096547607049,1.045467962003917,1.036214388809797,1.035175153515304,1.041935950315076,0.9944092447426414,1.048490414422779,1.015733011252985)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $vals[$j]
    }
}

Write-Host "Done updating vm_pu values"
